# Complete the "bai_tap" (homework) worksheet: convert the two essay-type
# rows (5 and 6) into multiple-choice rows like the others, and give every
# row a 4th answer option of "Các đáp án trên đều sai" (all of the above
# are wrong).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bai_tap")

# Row 2 - Colosseum question: add 4th (wrong) answer option
$ws.Range("H2").Value = "Các đáp án trên đều sai"

# Row 3 - Little Mermaid statue question: add 4th (wrong) answer option
$ws.Range("H3").Value = "Các đáp án trên đều sai"

# Row 4 - Myanmar pagoda question: add 4th (wrong) answer option
$ws.Range("H4").Value = "Các đáp án trên đều sai"

# Row 5 - Burj Khalifa question: turn into multiple choice like the others
$ws.Range("B5").Value = "Trắc nghiệm"
$ws.Range("E5").Value = "Shwedagon"
$ws.Range("F5").Value = "WatPhone"
$ws.Range("G5").Value = "That wlang"
$ws.Range("H5").Value = "Các đáp án trên đều sai"
$ws.Range("I5").Value = "a"

# Row 6 - St. Basil's Cathedral question: turn into multiple choice like the others
$ws.Range("B6").Value = "Trắc nghiệm"
$ws.Range("E6").Value = "Shwedagon"
$ws.Range("F6").Value = "WatPhone"
$ws.Range("G6").Value = "That wlang"
$ws.Range("H6").Value = "Các đáp án trên đều sai"
$ws.Range("I6").Value = "c"

# Selection / view bookkeeping: bai_tap becomes the active sheet/tab,
# with the view scrolled back to the top and C10 selected.
$ws.Range("C10").Select()
$ws.Activate()
